$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.359.14'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.976.03'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -6.13%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.12'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -5.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '124.96'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -7.46%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.974.03'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.134'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.19%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.14'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.434'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.45'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.38%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.458.95'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.482.79'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.970.27'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -6.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.12'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -6.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '430.88'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -6.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.00'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -6.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.655'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -6.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.11'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -6.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.89'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.25'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -6.32%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.46'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -8.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.11'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.88'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -7.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.21'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -7.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.01'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -10.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0936'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -9.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.25'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.949'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -8.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.49'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '49.31'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0659'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.22%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0356'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -8.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.74'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.42%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '372.93'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.18%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.106'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.662.61'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.42'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -7.78%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.234'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -6.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.71'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.93'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.95'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.106'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.22'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -8.51%  '
